$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.086.05"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.054.50"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'248.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").Value = "'0.656"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'55.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +17.45%  "
$ws.Range("D9").Value = "'61.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'0.379"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("D13").Value = "'15.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.40%  "
$ws.Range("D14").Value = "2.351.55"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "'0.817"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "'5.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "2.054.17"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "37.037.08"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'72.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  +8.25%  "
$ws.Range("D21").Value = "'14.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.44%  "
$ws.Range("D22").Value = "'5.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "'236.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").Value = "'169.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").Value = "'9.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "'20.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.48%  "
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "'4.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  +10.30%  "
$ws.Range("D33").Value = "'0.0624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").Value = "'4.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'0.0859"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.59%  "
$ws.Range("D37").Value = "'2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("D39").Value = "'1.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").Value = "'0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +22.61%  "
$ws.Range("D41").Value = "'18.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.52%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "'1.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").Value = "'95.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +42.61%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'14.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -51.16%  "
$ws.Range("D48").Value = "'2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.32%  "
$ws.Range("D49").Value = "1.296.04"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "'6.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.14%  "
